$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Save" header cell in H1
$ws.Range("H1").Value = "Save"

# Match the formatting of the neighboring header cell (G1: bold, bordered, centered)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Fill in the Save column values for rows 2-7
$saveValues = @(0, 0, 1, 0, 0, 1)
for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
